$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Julio de 2020 a las 06:12"

# Row 15 - Pakistan
$ws.Range("B15").Value = 246351
$ws.Range("C15").Value = 2752
$ws.Range("D15").Value = 153134
$ws.Range("E15").Value = 88094
$ws.Range("G15").Value = 65
$ws.Range("H15").Value = 5123

# Row 19 - Alemania
$ws.Range("D19").Value = 184500
$ws.Range("E19").Value = 5958

# Row 55 - Honduras
$ws.Range("B55").Value = 27053
$ws.Range("C55").Value = 669
$ws.Range("D55").Value = 2850
$ws.Range("E55").Value = 23453
$ws.Range("G55").Value = 46
$ws.Range("H55").Value = 750

# Rows 73-75: Kirguistan moves above Australia and Kenia because its
# total case count has grown past both. Australia and Kenia shift down
# one row with their figures unchanged, and Kirguistan's new totals land
# on row 73.
$ws.Range("A73").Value = "Kirguistan"
$ws.Range("B73").Value = 9672
$ws.Range("C73").Value = 314
$ws.Range("D73").Value = 3235
$ws.Range("E73").Value = 6312
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = 125

$ws.Range("A74").Value = "Australia"
$ws.Range("B74").Value = 9549
$ws.Range("C74").Value = 190
$ws.Range("D74").Value = 7730
$ws.Range("E74").Value = 1712
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 107

$ws.Range("A75").Value = "Kenia"
$ws.Range("B75").Value = 9448
$ws.Range("D75").Value = 2733
$ws.Range("E75").Value = 6534
$ws.Range("H75").Value = 181

# Row 87 - Haiti
$ws.Range("B87").Value = 6617
$ws.Range("C87").Value = 35
$ws.Range("D87").Value = 2590
$ws.Range("E87").Value = 3892
$ws.Range("G87").Value = 5
$ws.Range("H87").Value = 135

# Row 161 - Vietnam
$ws.Range("B161").Value = 370
$ws.Range("C161").Value = 1
$ws.Range("E161").Value = 20

# Row 170 - Mongolia
$ws.Range("D170").Value = 200
$ws.Range("E170").Value = 27
